$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.093714853270462584
$ws.Range("B1").Value = 0.093531712790372978
$ws.Range("A2").Value = -0.084097524412325875
$ws.Range("B2").Value = 0.083683130729807154
$ws.Range("A3").Value = -0.033977177575641448
$ws.Range("B3").Value = 0.033951005705759485
$ws.Range("A4").Value = -0.02595100583330634
$ws.Range("B4").Value = 0.025828150869527278
$ws.Range("A5").Value = -0.022828150929305124
$ws.Range("B5").Value = 0.022429891428383364
$ws.Range("A6").Value = -0.0086760630258648774
$ws.Range("B6").Value = 0.0085270224705027431
$ws.Range("A7").Value = 0.0014729773645716371
$ws.Range("B7").Value = -0.0014972255512262223
$ws.Range("A8").Value = 0.011497225386481347
$ws.Range("B8").Value = -0.011523252480719393
$ws.Range("A9").Value = 0.013523252431637101
$ws.Range("B9").Value = -0.013541610853327057
$ws.Range("A10").Value = 0.015541610807703776
$ws.Range("B10").Value = -0.01554100317556717
$ws.Range("A11").Value = 0.018541003115952748
$ws.Range("B11").Value = -0.018541766544696081
$ws.Range("A12").Value = 0.022041766479957925
$ws.Range("B12").Value = -0.022078446983656352
$ws.Range("A13").Value = 0.017846492813091253
$ws.Range("B13").Value = -0.017986732195576671
$ws.Range("A14").Value = 0.025986732078767893
$ws.Range("B14").Value = -0.02606369135902753
$ws.Range("A15").Value = -0.0080494065617493504
$ws.Range("B15").Value = 0.0080324206535324905
$ws.Range("A16").Value = -0.0060324206831596783
$ws.Range("B16").Value = 0.0060026646607012069
$ws.Range("A17").Value = -0.0040026646920985343
$ws.Range("B17").Value = 0.0039999999400555097
$ws.Range("A18").Value = -0.0034341884640234355
$ws.Range("B18").Value = 0.003360626855187121
$ws.Range("A19").Value = 0.0006393730867793046
$ws.Range("B19").Value = -0.0012011640716531602
$ws.Range("A20").Value = 0.005201164012420989
$ws.Range("B20").Value = -0.0053562533537920842
$ws.Range("A21").Value = 0.0093562532948716637
$ws.Range("B21").Value = -0.0095951682939476868
$ws.Range("A22").Value = -0.04570595322301152
$ws.Range("B22").Value = 0.045495200521422419
$ws.Range("A23").Value = -0.040495200609847792
$ws.Range("B23").Value = 0.040097956956954839
$ws.Range("A24").Value = -0.020097957267005917
$ws.Range("B24").Value = 0.019999999685413883
$ws.Range("A25").Value = -0.011753828454326154
$ws.Range("B25").Value = 0.011753862112023228
$ws.Range("A26").Value = -0.0092538621711621971
$ws.Range("B26").Value = 0.0092535631601275981
$ws.Range("A27").Value = -0.0067535632196875639
$ws.Range("B27").Value = 0.0067423213653805547
$ws.Range("A28").Value = -0.0047423214201254282
$ws.Range("B28").Value = 0.0047402808146204123
$ws.Range("A29").Value = 0.0022597190591087113
$ws.Range("B29").Value = -0.0022595366654538651
$ws.Range("A30").Value = 0.062259535792228871
$ws.Range("B30").Value = -0.062542993638926347
$ws.Range("A31").Value = 0.046068838155784775
$ws.Range("B31").Value = -0.046145434256697371
$ws.Range("A32").Value = -0.0040007170905571598
$ws.Range("B32").Value = 0.0039999999357274163